# Apply Trade #12 update across the workbook.
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.22
$summary.Range("B6").Value = 12
$summary.Range("B9").Value = 25

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 12
$status.Range("G4").Value = 25

# --- New trade row data ---
$tradeRow = @(12, "2026-02-17", "15:14:38", "MarketMaking", "DOWN", 0.02, 0.02, "CLOSED", 0, 0, 99.87, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.14)

# --- All Trades sheet ---
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(13, 2).NumberFormat = "@"
$allTrades.Cells.Item(13, 3).NumberFormat = "@"
for ($i = 0; $i -lt $tradeRow.Length; $i++) {
    $allTrades.Cells.Item(13, $i + 1).Value = $tradeRow[$i]
}

# --- MarketMaking sheet ---
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Cells.Item(13, 2).NumberFormat = "@"
$marketMaking.Cells.Item(13, 3).NumberFormat = "@"
for ($i = 0; $i -lt $tradeRow.Length; $i++) {
    $marketMaking.Cells.Item(13, $i + 1).Value = $tradeRow[$i]
}
